# Applies the "Updated cryptos list" data refresh described by the commit diff.
# Most rows only get refreshed Price (D) / Volume(1h) (E) values.
# Rows 42-47 also get re-ranked: "Aave" jumps to the top of that block (row 42)
# and RenderToken/Mantle/EnergySwap/VeChain/Hedera each shift down one row,
# so their Coin (B) and Link (C) values move together with new D/E figures.
#
# Numeric-looking Price values (e.g. "1.00", "0.997", "20.31") are written with a
# leading apostrophe so Excel keeps storing them as text instead of silently
# converting them to numbers (which would also round-trip "1.00" into "1").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.728.71'
$ws.Range("E2").Value = '  -0.23%  '

$ws.Range("D3").Value = '2.641.01'
$ws.Range("E3").Value = '  +1.22%  '

$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").Value = '''581.65'
$ws.Range("E5").Value = '  +0.19%  '

$ws.Range("D6").Value = '''143.65'
$ws.Range("E6").Value = '  -0.32%  '

$ws.Range("D8").Value = '''0.598'
$ws.Range("E8").Value = '  -0.25%  '

$ws.Range("E10").Value = '  +0.34%  '

$ws.Range("D11").Value = '''0.381'
$ws.Range("E11").Value = '  +2.55%  '

$ws.Range("D12").Value = '''0.157'
$ws.Range("E12").Value = '  -1.29%  '

$ws.Range("D13").Value = '3.110.34'
$ws.Range("E13").Value = '  +0.89%  '

$ws.Range("D14").Value = '''26.25'
$ws.Range("E14").Value = '  +11.06%  '

$ws.Range("D15").Value = '60.707.87'
$ws.Range("E15").Value = '  -0.25%  '

$ws.Range("E16").Value = '  +0.41%  '

$ws.Range("D17").Value = '2.648.81'
$ws.Range("E17").Value = '  +1.06%  '

$ws.Range("D18").Value = '''11.55'
$ws.Range("E18").Value = '  +1.86%  '

$ws.Range("D19").Value = '''4.73'
$ws.Range("E19").Value = '  +0.74%  '

$ws.Range("D20").Value = '''350.46'
$ws.Range("E20").Value = '  +0.13%  '

$ws.Range("D21").Value = '''6.87'
$ws.Range("E21").Value = '  -1.03%  '

$ws.Range("D22").Value = '''1.00'
$ws.Range("E22").Value = '  +0.43%  '

$ws.Range("D23").Value = '''0.527'
$ws.Range("E23").Value = '  +1.35%  '

$ws.Range("D24").Value = '''63.95'
$ws.Range("E24").Value = '  +1.05%  '

$ws.Range("D25").Value = '''0.997'
$ws.Range("E25").Value = '  +0.01%  '

$ws.Range("E26").Value = '  +0.96%  '

$ws.Range("D27").Value = '''8.29'
$ws.Range("E27").Value = '  +4.97%  '

$ws.Range("D28").Value = '''1.97'
$ws.Range("E28").Value = '  +7.85%  '

$ws.Range("D29").Value = '0.0₃0805'
$ws.Range("E29").Value = '  +0.44%  '

$ws.Range("D30").Value = '''6.81'
$ws.Range("E30").Value = '  +6.85%  '

$ws.Range("E31").Value = '  +0.06%  '

$ws.Range("D32").Value = '''165.05'
$ws.Range("E32").Value = '  +0.86%  '

$ws.Range("D33").Value = '''19.93'
$ws.Range("E33").Value = '  +1.88%  '

$ws.Range("D34").Value = '''4.53'
$ws.Range("E34").Value = '  +5.90%  '

$ws.Range("D35").Value = '''1.06'
$ws.Range("E35").Value = '  +3.67%  '

$ws.Range("D36").Value = '''1.32'
$ws.Range("E36").Value = '  +6.85%  '

$ws.Range("D37").Value = '''1.67'
$ws.Range("E37").Value = '  +2.08%  '

$ws.Range("D38").Value = '''337.31'
$ws.Range("E38").Value = '  +8.36%  '

$ws.Range("D39").Value = '''4.09'
$ws.Range("E39").Value = '  +4.41%  '

$ws.Range("D40").Value = '''0.906'
$ws.Range("E40").Value = '  +7.13%  '

$ws.Range("D41").Value = '''38.46'
$ws.Range("E41").Value = '  +1.28%  '

$ws.Range("B42").Value = 'Aave'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D42").Value = '''137.64'
$ws.Range("E42").Value = '  +1.85%  '

$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").Value = '''5.22'
$ws.Range("E43").Value = '  +3.57%  '

$ws.Range("B44").Value = 'Mantle'
$ws.Range("C44").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D44").Value = '''0.623'
$ws.Range("E44").Value = '  +2.48%  '

$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = '''20.31'
$ws.Range("E45").Value = '  +2.03%  '

$ws.Range("B46").Value = 'VeChain'
$ws.Range("C46").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D46").Value = '''0.0249'
$ws.Range("E46").Value = '  +2.67%  '

$ws.Range("B47").Value = 'Hedera'
$ws.Range("C47").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D47").Value = '''0.0562'
$ws.Range("E47").Value = '  +1.62%  '

$ws.Range("D48").Value = '''0.0997'
$ws.Range("E48").Value = '  +1.14%  '

$ws.Range("D49").Value = '''20.55'
$ws.Range("E49").Value = '  +0.84%  '

$ws.Range("D50").Value = '''0.999'
$ws.Range("E50").Value = '  +0.40%  '

$ws.Range("D51").Value = '2.092.29'
$ws.Range("E51").Value = '  +2.32%  '
